# ---------------------------------------------------------------------------
# Results.xlsx edit: correct the "Unoptimized" run-time measurements
# (row 3/4 of B:I) and add a second, derived "share of total runtime [%]"
# table (columns K:R) with a matching stacked percentage chart.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------------
# 1) Fix the measured values in the existing table (B3:I4)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 2191
$ws.Range("C3").Value = 873
$ws.Range("D3").Value = 4650
$ws.Range("E3").Value = 2096
$ws.Range("F3").Value = 178
$ws.Range("G3").Value = 64
$ws.Range("H3").Value = 2070
$ws.Range("I3").Value = 736

$ws.Range("B4").Value = 5.3

# ---------------------------------------------------------------------------
# 2) Mirror the header layout (A1:I2) into J1:R2 for the new derived table
# ---------------------------------------------------------------------------
$ws.Range("J3").Value = "Extraction"
$ws.Range("J4").Value = "Processing"
$ws.Range("J5").Value = "Classification"

# A1/A2 carry the plain "no border" style used for the spacer column (J);
# copy single (unmerged) cells so Excel does not invent extra half-border
# styles the way it would for a merged-range copy.
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Row 6: per-column totals used to normalise each measurement to a share
#    of the overall run time.
# ---------------------------------------------------------------------------
$ws.Range("K6").Formula = "=SUM(B3:B5)"
$ws.Range("L6").Formula = "=SUM(C3:C5)"
$ws.Range("M6").Formula = "=SUM(D3:D5)"
$ws.Range("N6").Formula = "=SUM(E3:E5)"
$ws.Range("O6").Formula = "=SUM(F3:F5)"
$ws.Range("P6").Formula = "=SUM(G3:G5)"
$ws.Range("Q6").Formula = "=SUM(H3:H5)"
$ws.Range("R6").Formula = "=SUM(I3:I5)"

# ---------------------------------------------------------------------------
# 4) K3:R5 -- each measurement expressed as a percentage of the column total
# ---------------------------------------------------------------------------
$ws.Range("K3").Formula = "=B3/K6*100"
$ws.Range("L3").Formula = "=C3/L6*100"
$ws.Range("M3").Formula = "=D3/M6*100"
$ws.Range("N3").Formula = "=E3/N6*100"
$ws.Range("O3").Formula = "=F3/O6*100"
$ws.Range("P3").Formula = "=G3/P6*100"
$ws.Range("Q3").Formula = "=H3/Q6*100"
$ws.Range("R3").Formula = "=I3/R6*100"

$ws.Range("K4").Formula = "=B4/K6*100"
$ws.Range("L4").Formula = "=C4/L6*100"
$ws.Range("M4").Formula = "=D4/M6*100"
$ws.Range("N4").Formula = "=E4/N6*100"
$ws.Range("O4").Formula = "=F4/O6*100"
$ws.Range("P4").Formula = "=G4/P6*100"
$ws.Range("Q4").Formula = "=H4/Q6*100"
$ws.Range("R4").Formula = "=I4/R6*100"

$ws.Range("K5").Formula = "=B5/K6*100"
$ws.Range("L5").Formula = "=C5/L6*100"
$ws.Range("M5").Formula = "=D5/M6*100"
$ws.Range("N5").Formula = "=E5/N6*100"
$ws.Range("O5").Formula = "=F5/O6*100"
$ws.Range("P5").Formula = "=G5/P6*100"
$ws.Range("Q5").Formula = "=H5/Q6*100"
$ws.Range("R5").Formula = "=I5/R6*100"

# ---------------------------------------------------------------------------
# 5) Copy cell formatting from the source table onto the new one so both
#    look the same (border + centred header row, plain data rows).
#    B2:I5 is a plain, unmerged block and can be copied in one shot; the
#    header row (B1:I1) is built by merging first and then applying the
#    border/alignment to the already-merged range -- doing it the other way
#    round makes Excel synthesise extra "half border" styles for the
#    internal merge edges.
# ---------------------------------------------------------------------------
$ws.Range("B2:I5").Copy()
$ws.Range("K2").PasteSpecial(-4122)

$ws.Range("K1:L1").Merge()
$ws.Range("M1:N1").Merge()
$ws.Range("O1:P1").Merge()
$ws.Range("Q1:R1").Merge()

$ws.Range("K1:R1").Borders.LineStyle = 1
$ws.Range("K1:R1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 6) Sheet-level bookkeeping to match the edited workbook
# ---------------------------------------------------------------------------
$ws.Range("B5").Select()

# ---------------------------------------------------------------------------
# 7) Add the second chart: stacked column chart of the percentage table
#    (K1:R5), mirroring "Diagramm 1" but showing the share of runtime.
# ---------------------------------------------------------------------------
$co2 = $ws.ChartObjects().Add(660, 142, 2000, 1040)
$chart2 = $co2.Chart
$chart2.ChartType = 52

$s1 = $chart2.SeriesCollection().NewSeries()
$s1.Name = "Extraction"
$s1.Values = $ws.Range("K3:R3")
$s1.XValues = $ws.Range("K1:R2")

$s2 = $chart2.SeriesCollection().NewSeries()
$s2.Name = "Processing"
$s2.Values = $ws.Range("K4:R4")
$s2.XValues = $ws.Range("K1:R2")

$s3 = $chart2.SeriesCollection().NewSeries()
$s3.Name = "Classification"
$s3.Values = $ws.Range("K5:R5")
$s3.XValues = $ws.Range("K1:R2")

$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Laufzeitmessungen auf dem ARM ohne aktive Optimierung "

$chart2.HasLegend = $true
$chart2.Legend.Position = -4107

$chart2.Axes(2).HasTitle = $true
$chart2.Axes(2).AxisTitle.Text = "Anteil in [%]"
$chart2.Axes(2).MaximumScale = 100

Write-Host "edit complete"
